$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Remove the "current" status (column H) from the pre-15 series rows that
#    are being superseded -- these rows keep everything else, only the
#    status ("current") marker in column H is cleared.
# ---------------------------------------------------------------------------
$rowsToClear = @(34, 35, 36, 58, 60, 61, 62, 63, 73, 74, 75, 97, 99, 100, 101)
foreach ($r in $rowsToClear) {
    $ws.Range("H$r").ClearContents()
}

# ---------------------------------------------------------------------------
# 2. Add a new row (103) anticipating 2050_TM152_FBP_PlusCrossing_15b.
#    It mirrors row 102 (same RTP2021 / FinalBlueprint / Plus / BAUS v2.25
#    / run182 / current pattern) except for the run-id in column C.
# ---------------------------------------------------------------------------
foreach ($col in @("A", "B", "C", "D", "E", "F", "G", "H")) {
    $ws.Range($col + "102").Copy($ws.Range($col + "103"))
}
$ws.Range("C103").Value = "2050_TM152_FBP_PlusCrossing_15b"

# ---------------------------------------------------------------------------
# 3. Update the saved view state to match (selection on the last entry row).
# ---------------------------------------------------------------------------
$ws.Range("C97").Select() | Out-Null
